$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.989.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.096.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.088.79"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -7.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.89"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.606.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.011.56"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.101.88"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "490.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.37"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.33"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.91"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.64"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.19%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.09"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.98%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "496.48"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.251.21"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0803"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0540"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.03%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.86"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.74%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.95%  "
